$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I11").Value = "sd"
$ws.Range("J11").Value = "Statement-non-opinion"
$ws.Range("I36").Value = "sd"
$ws.Range("J36").Value = "Statement-non-opinion"
$ws.Range("I44").Value = "sd"
$ws.Range("J44").Value = "Statement-non-opinion"
$ws.Range("I74").Value = "sd"
$ws.Range("J74").Value = "Statement-non-opinion"
$ws.Range("I75").Value = "sd"
$ws.Range("J75").Value = "Statement-non-opinion"
$ws.Range("I78").Value = "sd"
$ws.Range("J78").Value = "Statement-non-opinion"
$ws.Range("I83").Value = "sd"
$ws.Range("J83").Value = "Statement-non-opinion"
$ws.Range("I84").Value = "sv"
$ws.Range("J84").Value = "Statement-opinion"
$ws.Range("I101").Value = "%"
$ws.Range("J101").Value = "Uninterpretable"
$ws.Range("I103").Value = "sv"
$ws.Range("J103").Value = "Statement-opinion"
$ws.Range("I120").Value = "sv"
$ws.Range("J120").Value = "Statement-opinion"
$ws.Range("I124").Value = "sv"
$ws.Range("J124").Value = "Statement-opinion"
$ws.Range("I142").Value = "sd"
$ws.Range("J142").Value = "Statement-non-opinion"
$ws.Range("I160").Value = "sd"
$ws.Range("J160").Value = "Statement-non-opinion"
$ws.Range("I161").Value = "aa"
$ws.Range("J161").Value = "Agree/Accept"
$ws.Range("I162").Value = "sd"
$ws.Range("J162").Value = "Statement-non-opinion"
$ws.Range("I175").Value = "aa"
$ws.Range("J175").Value = "Agree/Accept"
$ws.Range("I179").Value = "sv"
$ws.Range("J179").Value = "Statement-opinion"
$ws.Range("I180").Value = "aa"
$ws.Range("J180").Value = "Agree/Accept"
$ws.Range("I186").Value = "aa"
$ws.Range("J186").Value = "Agree/Accept"
$ws.Range("I190").Value = "%"
$ws.Range("J190").Value = "Uninterpretable"
$ws.Range("I221").Value = "aa"
$ws.Range("J221").Value = "Agree/Accept"
$ws.Range("I222").Value = "%"
$ws.Range("J222").Value = "Uninterpretable"
$ws.Range("I285").Value = "sd"
$ws.Range("J285").Value = "Statement-non-opinion"
$ws.Range("I286").Value = "sv"
$ws.Range("J286").Value = "Statement-opinion"
$ws.Range("I299").Value = "sd"
$ws.Range("J299").Value = "Statement-non-opinion"
$ws.Range("I303").Value = "%"
$ws.Range("J303").Value = "Uninterpretable"
$ws.Range("I305").Value = "sv"
$ws.Range("J305").Value = "Statement-opinion"
$ws.Range("I313").Value = "sv"
$ws.Range("J313").Value = "Statement-opinion"
$ws.Range("I327").Value = "aa"
$ws.Range("J327").Value = "Agree/Accept"
$ws.Range("I375").Value = "sv"
$ws.Range("J375").Value = "Statement-opinion"
$ws.Range("I390").Value = "aa"
$ws.Range("J390").Value = "Agree/Accept"
$ws.Range("I397").Value = "sv"
$ws.Range("J397").Value = "Statement-opinion"
$ws.Range("I408").Value = "ba"
$ws.Range("J408").Value = "Appreciation"
$ws.Range("I414").Value = "sv"
$ws.Range("J414").Value = "Statement-opinion"
$ws.Range("I422").Value = "b"
$ws.Range("J422").Value = "Acknowledge (Backchannel)"
$ws.Range("I429").Value = "b"
$ws.Range("J429").Value = "Acknowledge (Backchannel)"
$ws.Range("I432").Value = "sv"
$ws.Range("J432").Value = "Statement-opinion"
$ws.Range("I444").Value = "%"
$ws.Range("J444").Value = "Uninterpretable"
$ws.Range("I455").Value = "aa"
$ws.Range("J455").Value = "Agree/Accept"
$ws.Range("I471").Value = "sd"
$ws.Range("J471").Value = "Statement-non-opinion"
$ws.Range("I474").Value = "sd"
$ws.Range("J474").Value = "Statement-non-opinion"
$ws.Range("I483").Value = "sd"
$ws.Range("J483").Value = "Statement-non-opinion"
$ws.Range("I489").Value = "aa"
$ws.Range("J489").Value = "Agree/Accept"
$ws.Range("I505").Value = "sd"
$ws.Range("J505").Value = "Statement-non-opinion"
$ws.Range("I510").Value = "sd"
$ws.Range("J510").Value = "Statement-non-opinion"
$ws.Range("I512").Value = "sv"
$ws.Range("J512").Value = "Statement-opinion"
$ws.Range("I525").Value = "sd"
$ws.Range("J525").Value = "Statement-non-opinion"
$ws.Range("I529").Value = "aa"
$ws.Range("J529").Value = "Agree/Accept"
$ws.Range("I550").Value = "b"
$ws.Range("J550").Value = "Acknowledge (Backchannel)"
$ws.Range("I562").Value = "aa"
$ws.Range("J562").Value = "Agree/Accept"
$ws.Range("I565").Value = "ba"
$ws.Range("J565").Value = "Appreciation"
$ws.Range("I568").Value = "ba"
$ws.Range("J568").Value = "Appreciation"
$ws.Range("I571").Value = "ba"
$ws.Range("J571").Value = "Appreciation"
$ws.Range("I592").Value = "ba"
$ws.Range("J592").Value = "Appreciation"
$ws.Range("I601").Value = "sd"
$ws.Range("J601").Value = "Statement-non-opinion"
$ws.Range("I632").Value = "sv"
$ws.Range("J632").Value = "Statement-opinion"
